$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 72.79331566666666
$ws.Range("H2").Value = 218.379947
$ws.Range("I2").Value = 0.2828741606141505
$ws.Range("J2").Value = 0.2828741606141506
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 16.535604
$ws.Range("N2").Value = 49.606812
$ws.Range("O2").Value = 0.2120453146491552
$ws.Range("P2").Value = 0.2120453146491552
$ws.Range("Q2").Value = 1203.681441710996
$ws.Range("R2").Value = 10833.13297539896
$ws.Range("S2").Value = 0.05998214039354322
$ws.Range("T2").Value = 0.05998214039354324
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 72.79331566666666
$ws.Range("H3").Value = 218.379947
$ws.Range("I3").Value = 0.2828741606141505
$ws.Range("J3").Value = 0.2828741606141506
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 40.62063066666667
$ws.Range("N3").Value = 121.861892
$ws.Range("O3").Value = 0.5209011059384622
$ws.Range("P3").Value = 0.5209011059384622
$ws.Range("Q3").Value = 2956.910390697747
$ws.Range("R3").Value = 26612.19351627972
$ws.Range("S3").Value = 0.1473494631053252
$ws.Range("T3").Value = 0.1473494631053252
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 72.79331566666666
$ws.Range("H4").Value = 218.379947
$ws.Range("I4").Value = 0.2828741606141505
$ws.Range("J4").Value = 0.2828741606141506
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 20.825229
$ws.Range("N4").Value = 62.475687
$ws.Range("O4").Value = 0.2670535794123827
$ws.Range("P4").Value = 0.2670535794123827
$ws.Range("Q4").Value = 1515.937468427621
$ws.Range("R4").Value = 13643.43721584859
$ws.Range("S4").Value = 0.07554255711528214
$ws.Range("T4").Value = 0.07554255711528216
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 117.1700846666667
$ws.Range("H5").Value = 351.510254
$ws.Range("I5").Value = 0.4553218801152877
$ws.Range("J5").Value = 0.4553218801152878
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 16.535604
$ws.Range("N5").Value = 49.606812
$ws.Range("O5").Value = 0.2120453146491552
$ws.Range("P5").Value = 0.2120453146491552
$ws.Range("Q5").Value = 1937.478120694472
$ws.Range("R5").Value = 17437.30308625025
$ws.Range("S5").Value = 0.09654887133569111
$ws.Range("T5").Value = 0.09654887133569114
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 117.1700846666667
$ws.Range("H6").Value = 351.510254
$ws.Range("I6").Value = 0.4553218801152877
$ws.Range("J6").Value = 0.4553218801152878
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 40.62063066666667
$ws.Range("N6").Value = 121.861892
$ws.Range("O6").Value = 0.5209011059384622
$ws.Range("P6").Value = 0.5209011059384622
$ws.Range("Q6").Value = 4759.52273442673
$ws.Range("R6").Value = 42835.70460984058
$ws.Range("S6").Value = 0.2371776709100333
$ws.Range("T6").Value = 0.2371776709100333
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 117.1700846666667
$ws.Range("H7").Value = 351.510254
$ws.Range("I7").Value = 0.4553218801152877
$ws.Range("J7").Value = 0.4553218801152878
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 20.825229
$ws.Range("N7").Value = 62.475687
$ws.Range("O7").Value = 0.2670535794123827
$ws.Range("P7").Value = 0.2670535794123827
$ws.Range("Q7").Value = 2440.093845132722
$ws.Range("R7").Value = 21960.8446061945
$ws.Range("S7").Value = 0.1215953378695634
$ws.Range("T7").Value = 0.1215953378695634
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 67.37122333333333
$ws.Range("H8").Value = 202.11367
$ws.Range("I8").Value = 0.2618039592705617
$ws.Range("J8").Value = 0.2618039592705618
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 16.535604
$ws.Range("N8").Value = 49.606812
$ws.Range("O8").Value = 0.2120453146491552
$ws.Range("P8").Value = 0.2120453146491552
$ws.Range("Q8").Value = 1114.02387003556
$ws.Range("R8").Value = 10026.21483032004
$ws.Range("S8").Value = 0.05551430291992088
$ws.Range("T8").Value = 0.05551430291992089
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 67.37122333333333
$ws.Range("H9").Value = 202.11367
$ws.Range("I9").Value = 0.2618039592705617
$ws.Range("J9").Value = 0.2618039592705618
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 40.62063066666667
$ws.Range("N9").Value = 121.861892
$ws.Range("O9").Value = 0.5209011059384622
$ws.Range("P9").Value = 0.5209011059384622
$ws.Range("Q9").Value = 2736.661580584849
$ws.Range("R9").Value = 24629.95422526365
$ws.Range("S9").Value = 0.1363739719231037
$ws.Range("T9").Value = 0.1363739719231037
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 67.37122333333333
$ws.Range("H10").Value = 202.11367
$ws.Range("I10").Value = 0.2618039592705617
$ws.Range("J10").Value = 0.2618039592705618
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 20.825229
$ws.Range("N10").Value = 62.475687
$ws.Range("O10").Value = 0.2670535794123827
$ws.Range("P10").Value = 0.2670535794123827
$ws.Range("Q10").Value = 1403.02115392681
$ws.Range("R10").Value = 12627.19038534129
$ws.Range("S10").Value = 0.06991568442753715
$ws.Range("T10").Value = 0.06991568442753716
